# Adds a new "Version 1.3" section (rows 20-22) to the compliance tracker
# worksheet, mirroring the existing "Version 1.2" section pattern.
# Cell values are written in the same order the original author entered
# them in, so newly interned shared strings land at the same indices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 20: section header -------------------------------------------------
$ws.Range("A20").Value = "Version 1.3"

# ---- Column A / B first (Control ID / Control Name) -------------------------
$ws.Range("A21").Value = "IAW-001"
$ws.Range("B21").Value = "Designated admin workstation "
$ws.Range("A22").Value = "AD – 007"
$ws.Range("B22").Value = "Review of PowerShell change logs"

# ---- Column C (NIST CSF 2.0) -------------------------------------------------
$ws.Range("C21").Value = "PR.IR-0"
$ws.Range("C22").Value = "DE.CM-01"

# ---- Column D (ISO 27001:2022) ----------------------------------------------
$ws.Range("D21").Value = "A 5.15"
$ws.Range("D22").Value = "A 8.16"

# ---- Column E (HIPAA) --------------------------------------------------------
$ws.Range("E21").Value = "164.308(a)(4)(ii)(b)"
$ws.Range("E22").Value = "164.308 (a)(6)(ii)"

# ---- Column F (PCI DSS 4.0), only row 22 has one -----------------------------
$ws.Range("F22").Value = "A3.5.1"

# ---- Column G (Coverage) -----------------------------------------------------
$ws.Range("G21").Value = "Fully covered"
$ws.Range("G22").Value = "Fully covered"

# ---- Column I (Evidence) -----------------------------------------------------
$ws.Range("I21").Value = "Test restuls, configuration"
$ws.Range("I22").Value = "Test restuls, configuration"

# ---- Row heights match the rest of the table ---------------------------------
$ws.Rows.Item(20).RowHeight = 50
$ws.Rows.Item(21).RowHeight = 50
$ws.Rows.Item(22).RowHeight = 50

# ---- Column I styling matches the rest of the table (wrap text) -------------
$ws.Range("I21").WrapText = $true
$ws.Range("I22").WrapText = $true

# ---- Coverage formulas (column J) --------------------------------------------
$ws.Range("J21").Formula = '=((COUNTIF(C21:F21, "<>"&"Unknown") + IF(G21="Fully covered", 1, IF(G21="Partially covered", 0.5, 0)))/ 5) * 100'
$ws.Range("J22").Formula = '=((COUNTIF(C22:F22, "<>"&"Unknown") + IF(G22="Fully covered", 1, IF(G22="Partially covered", 0.5, 0)))/ 5) * 100'

# ---- Hyperlinks on column D (ISO 27001 Annex references) --------------------
$ws.Hyperlinks.Add($ws.Range("D21"), "https://www.isms.online/iso-27001/annex-a-2022/5-15-access-control-2022/", [Type]::Missing, [Type]::Missing, "A 5.15")
$ws.Range("D21").Style = $ws.Range("D19").Style

$ws.Hyperlinks.Add($ws.Range("D22"), "https://www.isms.online/iso-27001/annex-a/8-16-monitoring-activities-2022/", [Type]::Missing, [Type]::Missing, "A 8.16")
$ws.Range("D22").Style = $ws.Range("D19").Style

# ---- Selection / view state mirrors where the author ended up editing -------
$ws.Range("A20:J22").Select()
$excel.ActiveWindow.ScrollRow = 17
